# Fill in the first task-summary row (Stage/Task/estimate/hours),
# record the week number and author, and update the cumulative-total
# label with the running total, matching the author's weekly submission.
#
# NOTE: the order of these assignments matters, since it controls the
# order new entries are appended to xl/sharedStrings.xml.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TASK SUMMARY SHEET")

$ws.Range("A3").Value = "Project Build"
$ws.Range("A14").Value = "Cumulative Total:160"
$ws.Range("C1").Value = "Richard Dobson"
$ws.Range("B3").Value = "Implementation of csv output functionality"
$ws.Range("E1").Value = 8
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 20

$ws.Range("B3").Select()
